$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, positioned before "总计", by
#    duplicating the "2021-Q4" sheet (same column layout / styling) and
#    renaming the copy.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBeforeCopy = $wb.Worksheets.Item("总计")
$totalIndexBeforeCopy = $totalSheetBeforeCopy.Index
$q4Sheet.Copy($totalSheetBeforeCopy, $null)

# NOTE: after Copy(), worksheet object references can end up repointed at the
# freshly inserted sheet instead of the original one they were fetched for,
# so re-resolve sheets afterwards. The new copy is inserted right before
# "总计", i.e. at the index "总计" used to have, so grab it positionally
# (robust against whatever "(2)" / "(copy)" suffix naming convention is used).
$newSheet = $wb.Worksheets.Item($totalIndexBeforeCopy)
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2. Fill in the fund-level data for 2022-Q1.
#    Columns B:G hold text-like values (fund codes, names and numbers that
#    must stay formatted as text, e.g. "004250"), so force a text number
#    format before assigning them. Column H stays numeric.
# ---------------------------------------------------------------------------
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "004250"
$newSheet.Range("C2").Value = "银河量化优选混合"
$newSheet.Range("D2").Value = "0.39"
$newSheet.Range("E2").Value = "80.03"
$newSheet.Range("F2").Value = "1.67"
$newSheet.Range("G2").Value = "0.0065"
$newSheet.Range("H2").Value = 8

$newSheet.Range("B3").Value = "005126"
$newSheet.Range("C3").Value = "银河量化稳进混合"
$newSheet.Range("D3").Value = "0.10"
$newSheet.Range("E3").Value = "78.20"
$newSheet.Range("F3").Value = "1.90"
$newSheet.Range("G3").Value = "0.0019"
$newSheet.Range("H3").Value = 8

# Row 3 of column A did not exist in the copied template (only one data row),
# so create it and mirror the header-ish formatting used by A2 (bold, thin
# border, centered / top aligned) to keep visual consistency with the rest
# of the "index" column.
$a3 = $newSheet.Range("A3")
$a3.Value = 1
$a3.Font.Bold = $true
$a3.Borders.LineStyle = 1
$a3.HorizontalAlignment = -4108
$a3.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: add a new leading row for 2022-Q1
#    and push the previous rows (2021-Q4, 2021-Q2) down by one position.
#    Values are written directly (bottom row up) instead of doing a
#    structural row-insert, which keeps the existing per-cell styling of
#    the "index" column (A) intact for the pre-existing rows.
# ---------------------------------------------------------------------------

# Row 4 <- old row 3 (2021-Q2)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

# Row 3 <- old row 2 (2021-Q4)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0

# Row 2 <- new (2022-Q1)
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.01

# A4 is a brand-new cell (previously empty), give it the same "index column"
# look (bold, thin border, centered / top aligned) used by A2/A3.
$a4 = $totalSheet.Range("A4")
$a4.Font.Bold = $true
$a4.Borders.LineStyle = 1
$a4.HorizontalAlignment = -4108
$a4.VerticalAlignment = -4160

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
